# Update sports_name (column D) values to include competition year suffix,
# matching the text already present in the certificate description (column B).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D22").Value = "Basketball competition,2011"
$ws.Range("D28").Value = "Rowing competition,2011"
$ws.Range("D65").Value = "Karate competition,2011"
$ws.Range("D70").Value = "Boxing competition,2011"
$ws.Range("D93").Value = "Weightlifting competition,2012"
$ws.Range("D106").Value = "Boxing competition,2011"
$ws.Range("D131").Value = "Basketball competition,2011"
$ws.Range("D145").Value = "Badminton competition,2011"
$ws.Range("D158").Value = "Triathlon competition,2012"
$ws.Range("D168").Value = "Weightlifting competition,2011"
$ws.Range("D171").Value = "Weightlifting competition,2011"
$ws.Range("D189").Value = "Handball competition,2011"
$ws.Range("D190").Value = "Judo competition,2011"
$ws.Range("D242").Value = "Badminton competition,2011"
$ws.Range("D246").Value = "Hockey competition,2011"
$ws.Range("D248").Value = "Softball competition,2011"
$ws.Range("D255").Value = "Athletics competition,2011"
$ws.Range("D284").Value = "Triathlon competition,2012"
$ws.Range("D295").Value = "Volleyball competition,2011"
$ws.Range("D309").Value = "Taekwondo competition,2011"
$ws.Range("D318").Value = "Taekwondo competition,2012"
$ws.Range("D338").Value = "Volleyball competition,2011"
$ws.Range("D344").Value = "Weightlifting competition,2012"
$ws.Range("D348").Value = "Taekwondo competition,2011"
$ws.Range("D351").Value = "Rowing competition,2011"
$ws.Range("D367").Value = "Badminton competition,2011"
$ws.Range("D380").Value = "Tennis competition,2011"
$ws.Range("D408").Value = "Triathlon competition,2011"
$ws.Range("D433").Value = "Softball competition,2011"
$ws.Range("D448").Value = "Rowing competition,2011"
$ws.Range("D449").Value = "Judo competition,2011"
